$wb = $excel.ActiveWorkbook

# --- TPSEE is the first sheet ---
$tpsee = $wb.Worksheets.Item(1)
$zoom  = $wb.Worksheets.Item("Zoom")

# --- Insert a new worksheet named "Sheet1" right after "TPSEE" ---
$newSheet = $wb.Worksheets.Add($null, $tpsee)
$newSheet.Name = "Sheet1"

# Fill in the new sheet's data (US / Florida / Miami / address)
$newSheet.Range("A1").Value = "US"
$newSheet.Range("B1").Value = "Florida"
$newSheet.Range("C1").Value = "Miami"
$newSheet.Range("D1").Value = "The Little Beet, 19501 Biscayne Blvd, Treats Food Hall, Floor 3, 33180, +1 305-359-5808"

# Give D1 the same "Consolas 9pt grey" look already used elsewhere in the workbook
# (copy formatting from TPSEE!C2, which already carries that exact style)
$tpsee.Range("C2").Copy()
$newSheet.Range("D1").PasteSpecial(-4122)

# Widen column D on the new sheet to fit the long address text
$newSheet.Columns.Item(4).ColumnWidth = 89.83333333333333

# Match the print setup used by the other data sheets (A4 / portrait)
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Selection on the new sheet covers the header row A1:D1
$newSheet.Range("A1:D1").Select()

# --- Re-style / widen columns on TPSEE ---
# B2 and F2 pick up the Consolas style already used by C2:E2
$tpsee.Range("C2").Copy()
$tpsee.Range("B2").PasteSpecial(-4122)
$tpsee.Range("F2").PasteSpecial(-4122)

# C2:E2 go back to plain/default formatting
$tpsee.Range("C2:E2").ClearFormats()

# Column B and F get wider to fit the new content
$tpsee.Columns.Item(2).ColumnWidth = 17.5
$tpsee.Columns.Item(6).ColumnWidth = 89.83333333333333

# Re-activate TPSEE and move the selection to F2 (was D2)
$tpsee.Activate()
$tpsee.Range("F2").Select()
